# test style then wrap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("Z100:Z101")
$rng.Style = "Normal"
$rng.WrapText = $true
